$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts")

# Insert a new row at 14 (pushes old rows 14-37 down to 15-38) for the new
# "8a / R8 / Fuse" BOM line, mirroring formatting from the row above (13)
# the way Excel's native row-insert carries formatting down.
$ws.Rows.Item(14).Insert()
$ws.Range("B13:L13").Copy()
$ws.Range("B14:L14").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H14:I14").ClearFormats()

$ws.Range("E14").Value = "8a"
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = "R8"
$ws.Range("H14").Value = "507-1818-1-ND  "
$ws.Range("I14").Value = "0ZCM0010FF2G"
$ws.Range("J14").Value = "FUSE PTC 100MA "

# Print area now needs to cover the extra row
$ws.PageSetup.PrintArea = "`$B`$2:`$L`$39"

# Make "Parts" the active/selected sheet with the new selection, and make
# sure "Header" is no longer flagged as the selected tab.
$ws.Range("J15").Select()
$ws.Activate()
